$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.060.32"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "1.787.69"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  +0.13%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "227.57"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("E6").Value = "  +1.21%  "

$ws.Range("E7").Value = "  +0.01%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "31.43"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.17"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("E11").Value = "  -2.79%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "2.044.17"
$ws.Range("E13").Value = "  -2.38%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "11.53"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +12.11%  "

$ws.Range("D15").Value = "1.783.43"
$ws.Range("E15").Value = "  -2.58%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.635"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").Value = "34.072.18"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("E18").Value = "  -3.06%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "69.56"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "253.01"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("E22").Value = "  +0.29%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("E25").Value = "  -2.68%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.17"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("E27").Value = "  -1.61%  "

$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("E31").Value = "  -1.56%  "

$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("E33").Value = "  -1.14%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.60"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").Value = "1.450.97"
$ws.Range("E36").Value = "  -8.17%  "

$ws.Range("E37").Value = "  -0.82%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.628"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("E39").Value = "  -1.67%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "83.54"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.05%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.82"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("E43").Value = "  -2.15%  "

$ws.Range("E44").Value = "  -3.22%  "

$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").Value = "1.944.25"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("E49").Value = "  +0.13%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "11.89"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.43%  "

$ws.Range("E51").Value = "  -3.49%  "
